$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @{
    "C3" = 0
    "D3" = 0.03807553614380086
    "E3" = 0.001649131501769719
    "C4" = 0
    "D4" = 0.07421083397630317
    "E4" = 0.002087028032034687
    "B5" = 2470.566758096668
    "B6" = 2482.584729487203
    "C7" = 1
    "D7" = 0.01812079843763099
    "E7" = 0.0007123348598917967
    "D8" = 0.01719593543409337
    "E8" = 0.0006516127066592582
    "C9" = 3
    "D9" = 0.01169292715860255
    "E9" = 0.001148203268079208
    "D10" = 0.02626145343863209
    "E10" = 0.001020999307578782
    "C11" = 3
    "D11" = 0.003022396703061827
    "E11" = 0.0009764666271430518
    "D12" = 0.01093761998669994
    "E12" = 0.0006440978654129237
    "C13" = 3
    "D13" = 0.01543155789994548
    "E13" = 0.001178157863618592
    "C14" = 3
    "D14" = 0.008538734288540894
    "E14" = 0.0009896646210438237
    "D15" = 0.01785136717517294
    "E15" = 0.0009051469308598815
    "C16" = 5
    "D16" = 0.04060174722303102
    "E16" = 0.001615856980916849
    "C17" = 3
    "D17" = 0.01159582284168615
    "E17" = 0.0008606274765313941
    "D18" = 0.02665576337872482
    "E18" = 0.001064855733249467
    "B19" = 2749.880085527932
    "C19" = 3
    "D19" = 0.01449933027765855
    "E19" = 0.001616024427231293
    "C21" = 1
    "D21" = 0.04851787932698472
    "E21" = 0.00121145340602076
    "C23" = 2
    "D23" = 0.001325443346966772
    "E23" = 0.0001743078728063171
    "C30" = 1
    "D30" = 0.0295946834879278
    "E30" = 0.0007853605849685846
    "C33" = 2
    "D33" = 0.004259860495121163
    "E33" = 0.001419953498373721
    "C34" = 1
    "D34" = 0.02124980265850609
    "E34" = 0.0008424112900538539
    "C37" = 3
    "D37" = 0.005900999224678612
    "E37" = 0.0005619999261598677
    "D39" = 0.004701757321759495
    "E39" = 0.001158403977824803
    "C40" = 3
    "D40" = 0.006531558847255968
    "E40" = 0.00340432764160008
    "D42" = 0.004899217546188956
    "E42" = 0.001140913675139894
    "D43" = 0.002803557666307301
    "E43" = 0.001268276087139017
    "C44" = 3
    "D44" = 0.005948902094598789
    "E44" = 0.0006953262188492091
    "C45" = 3
    "D45" = 0.007395361199218036
    "E45" = 0.0006530968331776967
    "C49" = 3
    "D49" = 0.0112606585593418
    "E49" = 0.0007295111026816115
    "D51" = 0.006682619188758948
    "E51" = 0.000506707252974176
    "D52" = 0.0130488628990849
    "E52" = 0.0004610066028577759
    "D53" = 0.008578015564787749
    "E53" = 0.0003750479024411025
    "C56" = 0
    "D56" = 0.02454167604123628
    "E56" = 0.001761364309179637
    "C57" = 1
    "D57" = 0.02308149543432147
    "E57" = 0.0007372586188144812
    "C58" = 0
    "D58" = 0.008149115499999781
    "E58" = 0.00184396044160579
    "D60" = 0.01910140741942075
    "E60" = 0.000745097717852869
}

foreach ($addr in $cellUpdates.Keys) {
    $ws.Range($addr).Value = $cellUpdates[$addr]
}
